# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# to match output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row => new value for sheet "展览"
$sheet1Updates = @{
    2  = 204
    3  = 268
    4  = 72
    5  = 258
    7  = 86
    10 = 50
    11 = 35
    12 = 107
    13 = 2372
    14 = 59
    15 = 38
    16 = 523
    17 = 544
    20 = 46
    22 = 1860
    23 = 3991
    26 = 1182
    27 = 227
    28 = 2086
    32 = 111
    34 = 418
    36 = 691
    37 = 436
    38 = 413
}

# Row => new value for sheet "全部类型"
$sheet4Updates = @{
    2  = 204
    3  = 268
    4  = 72
    5  = 258
    7  = 86
    10 = 50
    11 = 35
    12 = 107
    13 = 2372
    14 = 59
    16 = 38
    17 = 523
    18 = 544
    21 = 46
    23 = 1860
    24 = 3991
    27 = 1182
    28 = 227
    29 = 2086
    33 = 111
    35 = 418
    37 = 691
    38 = 436
    39 = 413
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

$wb.Save()
